$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data (old row 1 -> row 2, etc.)
$ws.Rows("1:1").Insert()

# New header row
$ws.Range("B1").Value = "measured"
$ws.Range("C1").Value = "computed"
$ws.Range("D1").Value = "percent err"

# Re-assert the "measured" (B) formulas, cell by cell, so the newly shifted
# rows keep clean, valid formulas tied to their own row.
$ws.Range("B2").Formula = "=A2/3.141596"
$ws.Range("B3").Formula = "=A3/3.141596"
$ws.Range("B4").Formula = "=A4/3.141596"
$ws.Range("B5").Formula = "=A5/3.141596"
$ws.Range("B6").Formula = "=A6/3.141596"
$ws.Range("B7").Formula = "=A7/3.141596"

# New "computed" column - the depth values actually measured on occluded
# trees, used to sanity check the estimator. C3 gets the new smaller
# Arial Unicode MS font, then the same look is copied down onto C4:C7.
$ws.Range("C2").Value = 0

$ws.Range("C3").Value = 19.119
$fnt = $ws.Range("C3").Font
$fnt.Name = "Arial Unicode MS"
$fnt.Size = 10

$ws.Range("C4").Value = 15.936
$ws.Range("C5").Value = 3.87
$ws.Range("C6").Value = 43.494
$ws.Range("C7").Value = 149.08

$ws.Range("C3").Copy()
$ws.Range("C4:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "percent err" column comparing computed vs. measured values
$ws.Range("D2").Formula = "=(C2-B2)/B2"
$ws.Range("D3").Formula = "=(C3-B3)/B3"
$ws.Range("D4").Formula = "=(C4-B4)/B4"
$ws.Range("D5").Formula = "=(C5-B5)/B5"
$ws.Range("D6").Formula = "=(C6-B6)/B6"
$ws.Range("D7").Formula = "=(C7-B7)/B7"

# Match the row heights used for the new comparison rows
$ws.Rows("3:7").RowHeight = 17

# Final selection left by the author after the edit
$ws.Range("B9").Select()
